$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: belt length 5m -> 4m
$ws.Range("B2").Value = "4m"

# Row 3: bearing quantity x20 -> x28
$ws.Range("B3").Value = "x28"

# Row 4: new Notes cell mirrors the (now x28) bearing quantity note
$ws.Range("B4").Value = "x28"

# Row 5: T-Nut description narrowed to M5 only
$ws.Range("A5").Value = "20 Series T Nuts M5 T Slot Nuts"

# Row 6: Hex socket bolts/nuts sizes changed from M4 M5 M6 to M3 M5
$ws.Range("A6").Value = "M3 M5 Hex Socket Bolts / Nuts"

# Row 8: stepper motor description updated
$ws.Range("A8").Value = "NEMA17 stepper motors"

# Match the rest of the Notes column's plain "Normal" style (B8 / B9
# previously fell back to the worksheet default; bring them in line with
# B2-B4, the same way the other Notes cells are styled)
$ws.Range("B8:B9").Style = "Normal"

# Move the active selection to A13
$ws.Range("A13").Select()
